# feat: add 2022-Q1 data
#
# 1) Insert a new "2022-Q1" worksheet (fund-level breakdown) right before
#    the existing "总计" (totals) sheet.
# 2) Insert a new row at the top of "总计"'s data (after the header) with
#    the 2022-Q1 summary figures, shifting the older quarters down and
#    renumbering the leading index column.

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# 1) New "2022-Q1" sheet, inserted immediately before "总计"
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Add($total)
$q1.Name = "2022-Q1"

# Header row (bold + border + centered, same look as the other quarter sheets)
$headerRange = $q1.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Fund-level rows (row 2..10): A=index (styled like header), B/C text,
# D/E/F/G numeric-looking text (kept as text, not converted to numbers),
# H = numeric rank.
$q1Rows = @(
    @("007449", "兴全多维价值混合A", "28.03", "84.40", "2.61", "0.7316", 5),
    @("011466", "兴业医疗保健混合A", "7.10", "84.68", "4.60", "0.3266", 6),
    @("920002", "中金精选股票A", "3.40", "86.43", "7.70", "0.2618", 5),
    @("001179", "德邦大健康灵活配置混合", "4.34", "92.02", "4.43", "0.1923", 5),
    @("007450", "兴全多维价值混合C", "4.48", "84.40", "2.61", "0.1169", 5),
    @("011467", "兴业医疗保健混合C", "2.22", "84.68", "4.60", "0.1021", 6),
    @("001110", "中欧瑾泉灵活配置混合 - A", "7.39", "22.25", "1.33", "0.0983", 9),
    @("001111", "中欧瑾泉灵活配置混合 - C", "2.20", "22.25", "1.33", "0.0293", 9),
    @("920922", "中金精选股票C", "0.14", "86.43", "7.70", "0.0108", 5)
)

$rowNum = 2
for ($i = 0; $i -lt $q1Rows.Count; $i++) {
    $data = $q1Rows[$i]

    $idxCell = $q1.Cells.Item($rowNum, 1)
    $idxCell.Font.Bold = $true
    $idxCell.Borders.LineStyle = 1
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Value = $i

    $textRange = $q1.Range("B${rowNum}:G${rowNum}")
    $textRange.NumberFormat = "@"

    $q1.Cells.Item($rowNum, 2).Value = $data[0]
    $q1.Cells.Item($rowNum, 3).Value = $data[1]
    $q1.Cells.Item($rowNum, 4).Value = $data[2]
    $q1.Cells.Item($rowNum, 5).Value = $data[3]
    $q1.Cells.Item($rowNum, 6).Value = $data[4]
    $q1.Cells.Item($rowNum, 7).Value = $data[5]
    $q1.Cells.Item($rowNum, 8).Value = $data[6]

    $rowNum = $rowNum + 1
}

# ---------------------------------------------------------------------
# 2) "总计" sheet: add the 2022-Q1 row on top of the existing history
#    (re-fetch the sheet reference by name -- inserting the new sheet
#    above shifted its position, and a stale reference would now point
#    at the just-created "2022-Q1" sheet instead)
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

$totalRows = @(
    @("2022-Q1", 9, 1.87),
    @("2021-Q4", 10, 1.09),
    @("2021-Q3", 13, 2.67),
    @("2021-Q2", 28, 15.82),
    @("2021-Q1", 22, 17.99),
    @("2020-Q4", 17, 8.18)
)

$rowNum = 2
for ($i = 0; $i -lt $totalRows.Count; $i++) {
    $data = $totalRows[$i]

    $idxCell = $total.Cells.Item($rowNum, 1)
    $idxCell.Font.Bold = $true
    $idxCell.Borders.LineStyle = 1
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Value = $i

    $total.Cells.Item($rowNum, 2).Value = $data[0]
    $total.Cells.Item($rowNum, 3).Value = $data[1]
    $total.Cells.Item($rowNum, 4).Value = $data[2]

    $rowNum = $rowNum + 1
}
